$wb = $excel.ActiveWorkbook

# --- Update the conversion note text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.63 = 6071.0 pesos`n✅ 6071.0 pesos = 1.63 = 955.74 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the numeric rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 612.42
$wsTasas.Range("O10").Value = 3718
$wsTasas.Range("O12").Value = 587.2
